$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "MaxDD = 0.914"
$ws.Range("A9").Value = "NetProfit = 45285822.7"
$ws.Range("A10").Value = "SharpeRatio = 0.800"
$ws.Range("A11").Value = "AnnualizedReturn = 0.598"
